# Add a new day column (CK) for 2024/12/06 to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set up the new column's width to match the other data columns (12) ---
# ColumnWidth uses Excel's "characters" unit, which gets converted internally;
# 11.1 round-trips to a stored width of exactly 12, matching the other cols.
$ws.Range("CK1").ColumnWidth = 11.1

# --- Header cell CK1: date label "2024/12/06" stored as literal text ---
# Assigning a date-looking string straight to .Value/.Value2 makes Excel
# auto-convert it into a real date serial number, which is not what the
# source file does (all the date headers are plain text). Instead: clone the
# format from the existing date header (CJ1, style index 1), then build the
# text via a formula and immediately collapse it back down to a static
# value so it is stored as literal text, not a formula.
$ws.Range("CJ1").Copy()
$ws.Range("CK1").PasteSpecial(-4122)
$ws.Range("CK1").Formula = "=""2024/12/06"""
$ws.Range("CK1").Copy()
$ws.Range("CK1").PasteSpecial(-4163)

# --- Data cells CK2:CK53: per-machine payout % for 2024/12/06 ---
# Styles in this sheet are value-based highlight colors applied by hand:
#   style 1 -> plain, style 2 -> yellow fill, style 3 -> light-blue fill.
# Reuse one already-styled cell of each kind as a format donor so the new
# cells land on the very same style index instead of minting new ones.
$styleSrc = @{1="A2"; 2="D2"; 3="N2"}

$data = @(
    @{row=2; style=1; val=152.3},
    @{row=3; style=1; val=177.1},
    @{row=4; style=1; val=153.3},
    @{row=5; style=1; val=163},
    @{row=6; style=3; val=127.2},
    @{row=7; style=1; val=191.7},
    @{row=8; style=3; val=136.9},
    @{row=9; style=1; val=143.1},
    @{row=10; style=1; val=163.2},
    @{row=11; style=1; val=480.3},
    @{row=12; style=1; val=166.8},
    @{row=13; style=1; val=154.5},
    @{row=14; style=2; val=123.9},
    @{row=15; style=3; val=136.5},
    @{row=16; style=3; val=133.2},
    @{row=17; style=1; val=236.6},
    @{row=18; style=2; val=106.5},
    @{row=19; style=1; val=154.2},
    @{row=20; style=1; val=201},
    @{row=21; style=1; val=156.7},
    @{row=22; style=1; val=158.5},
    @{row=23; style=1; val=224.1},
    @{row=24; style=1; val=147.9},
    @{row=25; style=1; val=189.4},
    @{row=26; style=1; val=162.5},
    @{row=27; style=1; val=261.3},
    @{row=28; style=1; val=140.9},
    @{row=29; style=2; val=108.7},
    @{row=30; style=1; val=164.5},
    @{row=31; style=1; val=144.4},
    @{row=32; style=1; val=209.3},
    @{row=33; style=1; val=158.2},
    @{row=34; style=2; val=118.2},
    @{row=35; style=1; val=247},
    @{row=36; style=2; val=113},
    @{row=37; style=1; val=249},
    @{row=38; style=1; val=168.9},
    @{row=39; style=3; val=127.3},
    @{row=40; style=3; val=138.4},
    @{row=41; style=1; val=170.4},
    @{row=42; style=1; val=240.5},
    @{row=43; style=3; val=138.5},
    @{row=44; style=3; val=128.4},
    @{row=45; style=3; val=128.9},
    @{row=46; style=1; val=207.4},
    @{row=47; style=1; val=174.1},
    @{row=48; style=1; val=251.2},
    @{row=49; style=1; val=164.8},
    @{row=50; style=1; val=179.3},
    @{row=51; style=3; val=125},
    @{row=52; style=1; val=153.7},
    @{row=53; style=1; val=157.5}
)

foreach ($d in $data) {
    $src = $styleSrc[$d.style]
    $ws.Range($src).Copy()
    $dst = "CK" + $d.row
    $ws.Range($dst).PasteSpecial(-4122)
    $ws.Range($dst).Value2 = $d.val
}
